# "Generate Report for handback"
#
# The localization report workbook has an "Overview" sheet plus one sheet
# per target locale ("zh-cn", "de-de"). Each locale sheet lists, per source
# file, the handoff info (columns A-D) and handback info (columns E-I).
# A handback run just completed for the "f3911a4e...md" source file in both
# locales, so:
#   - Status (col B, both the summary row and the "f3911a4e" row) flips from
#     "Ready for handoff" to "Handed back: in sync with en-US"
#   - Latest Target File / Latest Handback File (cols E/F) get populated
#     with the (hyperlinked) handed-back file names
#   - Latest Handback DateTime (col G) gets the timestamp of the handback

function Get-HyperlinkAddress($ws, $addr) {
    $found = ""
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $found = $hl.Address
        }
    }
    return $found
}

function Add-MatchingHyperlink($ws, $cellAddr, $url, $text) {
    $ws.Range($cellAddr).Value = $text
    $ws.Hyperlinks.Add($ws.Range($cellAddr), $url, "", "", $text) | Out-Null
    # Match the look of the workbook's existing hyperlink cells (underlined,
    # cornflower-blue font) instead of whatever default the new hyperlink
    # style would otherwise pick up.
    $ws.Range($cellAddr).Font.Underline = $true
    $ws.Range($cellAddr).Font.Color = 15570276
}

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# zh-cn -------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhTargetUrl = Get-HyperlinkAddress $wsZh '$A$2'
$zhHandbackUrl = Get-HyperlinkAddress $wsZh '$C$2'

$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("B3").Value = $newStatus

Add-MatchingHyperlink $wsZh "E2" $zhTargetUrl "f3911a4e-8629-4cee-9cae-e90d0eecdda9.md"
Add-MatchingHyperlink $wsZh "F2" $zhHandbackUrl "f3911a4e-8629-4cee-9cae-e90d0eecdda9.2eb3e130c00a8243a2656243f0a2a7c40878f9a8.zh-cn.xlf"
$wsZh.Range("G2").Value = "2016-01-28 11:28:09"

Add-MatchingHyperlink $wsZh "E3" $zhTargetUrl "f3911a4e-8629-4cee-9cae-e90d0eecdda9.md"
Add-MatchingHyperlink $wsZh "F3" $zhHandbackUrl "f3911a4e-8629-4cee-9cae-e90d0eecdda9.2eb3e130c00a8243a2656243f0a2a7c40878f9a8.zh-cn.xlf"
$wsZh.Range("G3").Value = "2016-01-28 11:28:09"

# de-de ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deTargetUrl = Get-HyperlinkAddress $wsDe '$A$2'
$deHandbackUrl = Get-HyperlinkAddress $wsDe '$C$2'

$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("B3").Value = $newStatus

Add-MatchingHyperlink $wsDe "E2" $deTargetUrl "f3911a4e-8629-4cee-9cae-e90d0eecdda9.md"
Add-MatchingHyperlink $wsDe "F2" $deHandbackUrl "f3911a4e-8629-4cee-9cae-e90d0eecdda9.2eb3e130c00a8243a2656243f0a2a7c40878f9a8.de-de.xlf"
$wsDe.Range("G2").Value = "2016-01-28 11:28:32"

Add-MatchingHyperlink $wsDe "E3" $deTargetUrl "f3911a4e-8629-4cee-9cae-e90d0eecdda9.md"
Add-MatchingHyperlink $wsDe "F3" $deHandbackUrl "f3911a4e-8629-4cee-9cae-e90d0eecdda9.2eb3e130c00a8243a2656243f0a2a7c40878f9a8.de-de.xlf"
$wsDe.Range("G3").Value = "2016-01-28 11:28:32"
